$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,35
$row2[0,0] = 3.792825088367897
$row2[0,1] = 1.546691642182422
$row2[0,2] = 1.688688498587012
$row2[0,3] = 1.081516395373035
$row2[0,4] = 1.271113164724866
$row2[0,5] = 1.862179354179809
$row2[0,6] = 4.217308298119121
$row2[0,7] = 2.177799723649859
$row2[0,8] = 1.94967491048222
$row2[0,9] = 1.158644663925531
$row2[0,10] = 1.051738910727824
$row2[0,11] = 1.344899113284217
$row2[0,12] = 1.043849219676489
$row2[0,13] = 1.05317222826534
$row2[0,14] = 0.8505002423756096
$row2[0,15] = 0.9358936184630279
$row2[0,16] = 1.106075010245341
$row2[0,17] = 1.275889375499644
$row2[0,18] = 1.929619331149395
$row2[0,19] = 1.600206971442398
$row2[0,20] = 1.301856881026457
$row2[0,21] = 1.323910226206807
$row2[0,22] = 1.172471701078861
$row2[0,23] = 1.351735241609772
$row2[0,24] = 1.406662577814643
$row2[0,25] = 1.880528686459688
$row2[0,26] = 1.777253694197898
$row2[0,27] = 2.041188729605869
$row2[0,28] = 1.36513940518563
$row2[0,29] = 1.560409587992492
$row2[0,30] = 1.003201609179825
$row2[0,31] = 1.473103553624098
$row2[0,32] = 1.731918962280045
$row2[0,33] = 7.390936627678586
$row2[0,34] = 2.081786204327535
$ws.Range("B2:AJ2").Value = $row2

$row3 = New-Object 'object[,]' 1,35
$row3[0,0] = 5.57554810382968
$row3[0,1] = 0.9616568858773247
$row3[0,2] = 1.30671255003566
$row3[0,3] = 0.6689411156364747
$row3[0,4] = 0.8389173022727757
$row3[0,5] = 0.9154935107345173
$row3[0,6] = 4.140016614340619
$row3[0,7] = 1.301037697724458
$row3[0,8] = 1.233914165336793
$row3[0,9] = 0.5353524692443765
$row3[0,10] = 0.4831433054144341
$row3[0,11] = 0.5628273149101107
$row3[0,12] = 0.4147661438194007
$row3[0,13] = 0.5130690958340587
$row3[0,14] = 0.3158494005154383
$row3[0,15] = 0.5462863584759837
$row3[0,16] = 0.5903789050660904
$row3[0,17] = 0.6173073900047488
$row3[0,18] = 1.044443586791306
$row3[0,19] = 0.851910146715721
$row3[0,20] = 0.6388790428753343
$row3[0,21] = 0.5748645361840525
$row3[0,22] = 0.4961483590010816
$row3[0,23] = 0.6115925526840317
$row3[0,24] = 0.8800274962732321
$row3[0,25] = 1.219838824567253
$row3[0,26] = 0.9519792161616931
$row3[0,27] = 1.078894259070486
$row3[0,28] = 0.7474179547294174
$row3[0,29] = 0.8634578096904488
$row3[0,30] = 0.4890578176541391
$row3[0,31] = 1.000423185388284
$row3[0,32] = 1.298616981540407
$row3[0,33] = 9.850804531664439
$row3[0,34] = 1.212763345597955
$ws.Range("B3:AJ3").Value = $row3

$row4 = New-Object 'object[,]' 1,35
$row4[0,0] = 6.785091908209915
$row4[0,1] = 1.830652960054959
$row4[0,2] = 2.149978415893389
$row4[0,3] = 1.268585576322448
$row4[0,4] = 1.518935107483683
$row4[0,5] = 2.411321500357255
$row4[0,6] = 7.96816516688512
$row4[0,7] = 2.562365772017433
$row4[0,8] = 2.319518689017587
$row4[0,9] = 1.289748811020303
$row4[0,10] = 1.157403209735369
$row4[0,11] = 1.484237851259173
$row4[0,12] = 1.126065568381799
$row4[0,13] = 1.17149973943205
$row4[0,14] = 0.9072549289405555
$row4[0,15] = 1.08366307057903
$row4[0,16] = 1.25377397478025
$row4[0,17] = 1.427681710736513
$row4[0,18] = 2.194149759960548
$row4[0,19] = 1.841220233773073
$row4[0,20] = 1.452728115631185
$row4[0,21] = 1.443332089997678
$row4[0,22] = 1.273127285062347
$row4[0,23] = 1.483655490303729
$row4[0,24] = 1.659261282023139
$row4[0,25] = 2.288442110368584
$row4[0,26] = 2.020240742238044
$row4[0,27] = 2.350418343437572
$row4[0,28] = 1.559614109057333
$row4[0,29] = 1.78337810668807
$row4[0,30] = 1.116060490148101
$row4[0,31] = 1.780696669722945
$row4[0,32] = 2.164705374098356
$row4[0,33] = 13.14380415315603
$row4[0,34] = 2.42649567819505
$ws.Range("B4:AJ4").Value = $row4

$row5 = New-Object 'object[,]' 1,35
$row5[0,0] = 0.2363636363636364
$row5[0,1] = 0.7592592592592593
$row5[0,2] = 0.7222222222222222
$row5[0,3] = 0.9122807017543859
$row5[0,4] = 0.8421052631578947
$row5[0,5] = 0.75
$row5[0,6] = 0.5714285714285714
$row5[0,7] = 0.5192307692307693
$row5[0,8] = 0.5740740740740741
$row5[0,9] = 0.9
$row5[0,10] = 0.9642857142857143
$row5[0,11] = 0.8666666666666667
$row5[0,12] = 0.9444444444444444
$row5[0,13] = 0.9642857142857143
$row5[0,14] = 1
$row5[0,15] = 0.9285714285714286
$row5[0,16] = 0.9107142857142857
$row5[0,17] = 0.8269230769230769
$row5[0,18] = 0.625
$row5[0,19] = 0.7551020408163265
$row5[0,20] = 0.9090909090909091
$row5[0,21] = 0.8928571428571429
$row5[0,22] = 0.9464285714285714
$row5[0,23] = 0.8392857142857143
$row5[0,24] = 0.7857142857142857
$row5[0,25] = 0.673469387755102
$row5[0,26] = 0.6363636363636364
$row5[0,27] = 0.6666666666666666
$row5[0,28] = 0.8545454545454545
$row5[0,29] = 0.7321428571428571
$row5[0,30] = 0.9464285714285714
$row5[0,31] = 0.7678571428571429
$row5[0,32] = 0.7142857142857143
$row5[0,33] = 0.3478260869565217
$row5[0,34] = 0.5660377358490566
$ws.Range("B5:AJ5").Value = $row5

$row6 = New-Object 'object[,]' 1,35
$row6[0,0] = 0.3818181818181818
$row6[0,1] = 0.8703703703703703
$row6[0,2] = 0.8703703703703703
$row6[0,3] = 0.9298245614035088
$row6[0,4] = 0.9298245614035088
$row6[0,5] = 0.85
$row6[0,6] = 0.6190476190476191
$row6[0,7] = 0.6153846153846154
$row6[0,8] = 0.6666666666666666
$row6[0,9] = 0.9399999999999999
$row6[0,10] = 0.9821428571428571
$row6[0,11] = 0.9555555555555556
$row6[0,12] = 1
$row6[0,13] = 1
$row6[0,14] = 1
$row6[0,15] = 0.9821428571428571
$row6[0,16] = 0.9821428571428571
$row6[0,17] = 0.9423076923076923
$row6[0,18] = 0.7678571428571429
$row6[0,19] = 0.8571428571428571
$row6[0,20] = 0.9272727272727272
$row6[0,21] = 0.9642857142857143
$row6[0,22] = 0.9821428571428571
$row6[0,23] = 0.9285714285714286
$row6[0,24] = 0.8928571428571429
$row6[0,25] = 0.7551020408163265
$row6[0,26] = 0.8181818181818182
$row6[0,27] = 0.7708333333333334
$row6[0,28] = 0.8909090909090909
$row6[0,29] = 0.8571428571428571
$row6[0,30] = 0.9821428571428571
$row6[0,31] = 0.8392857142857143
$row6[0,32] = 0.8035714285714286
$row6[0,33] = 0.4347826086956522
$row6[0,34] = 0.6981132075471698
$ws.Range("B6:AJ6").Value = $row6

$row7 = New-Object 'object[,]' 1,35
$row7[0,0] = 0.6363636363636364
$row7[0,1] = 0.8888888888888888
$row7[0,2] = 0.9074074074074074
$row7[0,3] = 1
$row7[0,4] = 0.9473684210526315
$row7[0,5] = 0.9
$row7[0,6] = 0.7619047619047619
$row7[0,7] = 0.7115384615384616
$row7[0,8] = 0.7777777777777778
$row7[0,9] = 1
$row7[0,10] = 1
$row7[0,11] = 1
$row7[0,12] = 1
$row7[0,13] = 1
$row7[0,14] = 1
$row7[0,15] = 1
$row7[0,16] = 1
$row7[0,17] = 1
$row7[0,18] = 0.875
$row7[0,19] = 0.9387755102040817
$row7[0,20] = 0.9636363636363636
$row7[0,21] = 0.9821428571428571
$row7[0,22] = 1
$row7[0,23] = 1
$row7[0,24] = 0.9285714285714286
$row7[0,25] = 0.8367346938775511
$row7[0,26] = 0.8545454545454545
$row7[0,27] = 0.875
$row7[0,28] = 0.9454545454545454
$row7[0,29] = 0.9464285714285714
$row7[0,30] = 1
$row7[0,31] = 0.8928571428571429
$row7[0,32] = 0.8392857142857143
$row7[0,33] = 0.5
$row7[0,34] = 0.7169811320754716
$ws.Range("B7:AJ7").Value = $row7

$row8 = New-Object 'object[,]' 1,35
$row8[0,0] = 0.7818181818181819
$row8[0,1] = 0.9629629629629629
$row8[0,2] = 0.9814814814814815
$row8[0,3] = 1
$row8[0,4] = 0.9824561403508771
$row8[0,5] = 0.95
$row8[0,6] = 0.8095238095238095
$row8[0,7] = 0.9038461538461539
$row8[0,8] = 0.9444444444444444
$row8[0,9] = 1
$row8[0,10] = 1
$row8[0,11] = 1
$row8[0,12] = 1
$row8[0,13] = 1
$row8[0,14] = 1
$row8[0,15] = 1
$row8[0,16] = 1
$row8[0,17] = 1
$row8[0,18] = 0.9642857142857143
$row8[0,19] = 0.9795918367346939
$row8[0,20] = 1
$row8[0,21] = 1
$row8[0,22] = 1
$row8[0,23] = 1
$row8[0,24] = 1
$row8[0,25] = 0.9387755102040817
$row8[0,26] = 0.9636363636363636
$row8[0,27] = 0.9166666666666666
$row8[0,28] = 1
$row8[0,29] = 0.9821428571428571
$row8[0,30] = 1
$row8[0,31] = 0.9821428571428571
$row8[0,32] = 0.9285714285714286
$row8[0,33] = 0.6521739130434783
$row8[0,34] = 0.9245283018867925
$ws.Range("B8:AJ8").Value = $row8

Write-Output "done"